$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Società"
$ws.Range("F1").Value = "Via"
$ws.Range("I1").Value = "Stato/Provincia"
$ws.Range("J1").Value = "Email"
$ws.Range("N1").Value = "Partita Iva"
$ws.Range("O1").Value = "Codice Fiscale"
$ws.Range("P1").Value = "Nome"
$ws.Range("Q1").Value = "Cognome"
$ws.Range("R1").Value = "Forma Giuridica"
$ws.Range("S1").Value = "Fonte del lead"
$ws.Range("T1").Value = "Km Annui Percorsi"
$ws.Range("U1").Value = "Reddito annuale"
$ws.Range("V1").Value = "Consenso Privacy"
$ws.Range("W1").Value = "Consenso Promozioni e Newsletter"
$ws.Range("Y1").Value = "Data Richiesta"
